# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 533b14e5-b481-45f0-9339-8a423ddbac09 file row after a new
# handback report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     533b14e5... row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 13:00:58"

# --- zh-cn sheet: update Correspond Handoff/Handback datetimes for the
#     533b14e5... row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-31 13:00:52"
$wsZhCn.Range("K3").Value = "2016-08-31 13:01:37"

# --- de-de sheet: update Correspond Handoff/Handback datetimes for the
#     533b14e5... row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-31 13:00:58"
$wsDeDe.Range("K3").Value = "2016-08-31 13:01:45"
